$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force plain-text storage for Price cells whose new value would otherwise
# be auto-parsed as a number by Excel (matches the source data which keeps
# these as literal text, e.g. "1.00", "324.70").
$textCells = @("D4", "D5", "D6", "D8", "D10", "D11", "D13", "D14", "D17", "D19", "D22", "D23", "D28", "D29", "D30", "D32", "D33", "D35", "D37", "D38", "D39", "D40", "D43", "D44", "D47", "D48", "D49", "D50", "D51")
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = "47.349.17"
$ws.Range("E2").Value = "  +6.04%  "
$ws.Range("D3").Value = "2.509.18"
$ws.Range("E3").Value = "  +3.74%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "324.70"
$ws.Range("E5").Value = "  +2.58%  "
$ws.Range("D6").Value = "106.23"
$ws.Range("E6").Value = "  +5.01%  "
$ws.Range("E7").Value = "  +2.16%  "
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("E9").Value = "  +2.66%  "
$ws.Range("D10").Value = "36.76"
$ws.Range("E10").Value = "  +4.32%  "
$ws.Range("D11").Value = "0.0819"
$ws.Range("E11").Value = "  +2.56%  "
$ws.Range("E12").Value = "  +0.76%  "
$ws.Range("D13").Value = "18.40"
$ws.Range("E13").Value = "  -1.03%  "
$ws.Range("D14").Value = "7.18"
$ws.Range("E14").Value = "  +4.29%  "
$ws.Range("D15").Value = "2.903.16"
$ws.Range("E15").Value = "  +3.75%  "
$ws.Range("D16").Value = "2.537.77"
$ws.Range("E16").Value = "  +4.78%  "
$ws.Range("D17").Value = "0.846"
$ws.Range("E17").Value = "  +2.24%  "
$ws.Range("D18").Value = "47.270.92"
$ws.Range("E18").Value = "  +6.29%  "
$ws.Range("D19").Value = "12.80"
$ws.Range("E19").Value = "  +5.06%  "
$ws.Range("E20").Value = "  +2.71%  "
$ws.Range("D21").Value = "0.0₃0943"
$ws.Range("E21").Value = "  +2.85%  "
$ws.Range("D22").Value = "70.99"
$ws.Range("E22").Value = "  +3.58%  "
$ws.Range("D23").Value = "252.81"
$ws.Range("E23").Value = "  +4.59%  "
$ws.Range("E24").Value = "  +5.65%  "
$ws.Range("E26").Value = "  +4.90%  "
$ws.Range("E27").Value = "  -0.01%  "
$ws.Range("D28").Value = "9.98"
$ws.Range("E28").Value = "  +5.15%  "
$ws.Range("D29").Value = "2.22"
$ws.Range("E29").Value = "  -3.02%  "
$ws.Range("D30").Value = "35.21"
$ws.Range("E30").Value = "  +5.73%  "
$ws.Range("E31").Value = "  +9.02%  "
$ws.Range("D32").Value = "49.84"
$ws.Range("E32").Value = "  +3.30%  "
$ws.Range("D33").Value = "19.85"
$ws.Range("E34").Value = "  +3.52%  "
$ws.Range("D35").Value = "0.0779"
$ws.Range("E35").Value = "  +2.06%  "
$ws.Range("E36").Value = "  +0.14%  "
$ws.Range("B37").Value = "ARBITRUM"
$ws.Range("C37").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D37").Value = "1.95"
$ws.Range("E37").Value = "  +4.15%  "
$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D38").Value = "4.65"
$ws.Range("E38").Value = "  +4.93%  "
$ws.Range("D39").Value = "2.97"
$ws.Range("E39").Value = "  +4.60%  "
$ws.Range("D40").Value = "124.15"
$ws.Range("E40").Value = "  -1.19%  "
$ws.Range("E41").Value = "  +2.50%  "
$ws.Range("E42").Value = "  +2.58%  "
$ws.Range("D43").Value = "21.61"
$ws.Range("E43").Value = "  +3.43%  "
$ws.Range("D44").Value = "0.0298"
$ws.Range("E44").Value = "  +3.43%  "
$ws.Range("D45").Value = "1.980.37"
$ws.Range("E45").Value = "  +2.07%  "
$ws.Range("E46").Value = "  +3.64%  "
$ws.Range("D47").Value = "2.13"
$ws.Range("E47").Value = "  +0.26%  "
$ws.Range("D48").Value = "1.81"
$ws.Range("E48").Value = "  +4.06%  "
$ws.Range("B49").Value = "FraxShare"
$ws.Range("C49").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D49").Value = "9.09"
$ws.Range("E49").Value = "  -0.37%  "
$ws.Range("B50").Value = "THORChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D50").Value = "5.43"
$ws.Range("E50").Value = "  +18.18%  "
$ws.Range("D51").Value = "80.19"
$ws.Range("E51").Value = "  +6.91%  "
